# Fix elective course scheduling to use common time slots for both sections
$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("Section_A")
$wsB = $wb.Worksheets.Item("Section_B")

# --- Section_A updates ---
$wsA.Range("B2").Value = "CS261"
$wsA.Range("E2").Value = "Free"
$wsA.Range("F2").Value = "CS264"

$wsA.Range("B3").Value = "MA261"
$wsA.Range("C3").Value = "CS264"
$wsA.Range("D3").Value = "CS263"
$wsA.Range("E3").Value = "CS264"

$wsA.Range("C5").Value = "Free"
$wsA.Range("D5").Value = "Free"
$wsA.Range("E5").Value = "CS261"
$wsA.Range("F5").Value = "CS263"

$wsA.Range("C6").Value = "CS263"
$wsA.Range("E6").Value = "Free"
$wsA.Range("F6").Value = "CS264 (Tutorial)"

$wsA.Range("B7").Value = "Free"
$wsA.Range("C7").Value = "MA261"
$wsA.Range("F7").Value = "Free"

# --- Section_B updates ---
$wsB.Range("B2").Value = "Free"
$wsB.Range("C2").Value = "CS261"
$wsB.Range("E2").Value = "Free"

$wsB.Range("B3").Value = "Free"

$wsB.Range("B5").Value = "CS263"
$wsB.Range("C5").Value = "CS263"
$wsB.Range("D5").Value = "MA261"
$wsB.Range("E5").Value = "Free"

$wsB.Range("B6").Value = "MA261"
$wsB.Range("C6").Value = "Free"
$wsB.Range("D6").Value = "CS264 (Tutorial)"
$wsB.Range("E6").Value = "CS261"
$wsB.Range("F6").Value = "Free"

$wsB.Range("B7").Value = "Free"
$wsB.Range("C7").Value = "CS264"
$wsB.Range("D7").Value = "CS264"
$wsB.Range("E7").Value = "CS264"
$wsB.Range("F7").Value = "CS261"
